$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column header D1: Tobacco -> Cassava
$ws.Range("D1").Value = "Cassava"

# Insert a new row at position 3, pushing old row 3 (yield_per_ton) down to row 4
$ws.Rows.Item(3).Insert() | Out-Null

# Fill the newly inserted row 3: Revenue (R$/ton)
$ws.Range("A3").Value = "Revenue (R`$/ton)"
$ws.Range("B3").Value = 664
$ws.Range("C3").Value = 1845
$ws.Range("D3").Value = 440

# Old row3 (now row4) had "yield_per_ton" -> becomes "Cost (R$/ton)"
$ws.Range("A4").Value = "Cost (R`$/ton)"
$ws.Range("B4").Value = 448
$ws.Range("C4").Value = 1351
$ws.Range("D4").Value = 333

# Row2 "efficiency" -> "Yield (ton/ha)" updated last
$ws.Range("A2").Value = "Yield (ton/ha)"
$ws.Range("B2").Value = 0.764
$ws.Range("C2").Value = 14.428
$ws.Range("D2").Value = 11.392

$ws.Columns.Item(1).ColumnWidth = 16

$ws.Range("F13").Select() | Out-Null
